$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "Conc (ug/m3)"
$ws.Range("C1").Value = "Conc sci (ug/m3)"
